# Atualizacao de bases das ligas, do dia: 09-05-2024 as 19:13
# Swap the data (columns B:AB) between rows 11/12 and rows 83/84,
# keeping the index column A unchanged for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($ws, $row1, $row2) {
    $cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB")
    foreach ($col in $cols) {
        $addr1 = "$col$row1"
        $addr2 = "$col$row2"
        $v1 = $ws.Range($addr1).Value()
        $v2 = $ws.Range($addr2).Value()
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}

Swap-RowData $ws 11 12
Swap-RowData $ws 83 84
